# Update COVID-19 "paises" (countries) stats and re-sort by total cases.
#
# The source workbook lists one row per country with columns:
#   A = Pais (country)        E = Recuperados
#   B = Casos totales          F = Casos criticos
#   C = Nuevos casos           G = Muertes hoy
#   D = Casos activos          H = Muertes
# sorted descending by column B ("Casos totales"). Data rows span 4..216.
#
# This script updates the 9 countries whose figures changed in the new
# data pull, re-sorts the table (since a couple of the updates change
# relative rank), and refreshes the "last updated" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Push the new per-country figures -------------------------------
# Row numbers below are this country's row *before* the re-sort below;
# writing straight to B:H on that row is unambiguous because A (the
# country name) is left untouched.
$updates = @{
    4   = @(574138, 13838, 33754, 517424, 12547, 855, 22960)   # Estados Unidos
    7   = @(136779, 4188,  27718, 94094,  6821,  574, 14967)   # Francia
    8   = @(128092, 238,   64300, 60754,  4895,  16,  3038)    # Alemania
    16  = @(25546,  1163,  7659,  17120,  557,   50,  767)     # Canada
    25  = @(10453,  1248,  1181,  8914,   0,     27,  358)     # India
    31  = @(6633,   333,   914,   5388,   231,   15,  331)     # Rumania
    75  = @(1091,   140,   138,   941,    21,    2,   12)      # Kazajistan
    112 = @(274,    2,     5,     266,    7,     0,   3)       # Montenegro
    155 = @(62,     21,    2,     56,     0,     0,   4)       # Birmania
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # column B is index 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# --- 2. Re-sort the data rows by "Casos totales" (column B), descending -
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)

# --- 3. Refresh the "last updated" banner --------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 19:52"
